# Commit: implemented the new TestEnvironment and TestNgTestClassTemplate framework
#
# The underlying data provider workbook was renamed/repurposed from the
# "AddNewTitle" test scenario to the "AddNewDept" test scenario:
#   - the worksheet itself is renamed to match the new test name
#   - the last active selection left in the sheet view moves from A2 to C28
#     (the cell the author was last looking at when the workbook was saved)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: TestAddNewTitle -> TestAddNewDept
$ws.Name = "TestAddNewDept"

# Update the saved selection/active cell: A2 -> C28
$ws.Range("C28").Select()
